# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scrape values (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column prices are plain text in the source data (note values like
# "1.000" / "30.489.79" that must keep their literal digits/dots), so each
# D cell is switched to the Text number format before the value is written;
# otherwise Excel would auto-coerce the numeric-looking strings into numbers
# and silently drop formatting (e.g. "1.000" -> 1).
# E-column percentages already contain padding spaces and a trailing "%" so
# Excel keeps them as text natively; no format coercion is required there.

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '30.489.79'
$ws.Cells.Item(2, 5).Value = '  +0.26%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.914.60'
$ws.Cells.Item(3, 5).Value = '  -0.13%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.9997'
$ws.Cells.Item(4, 5).Value = '  +0.02%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '244.84'
$ws.Cells.Item(5, 5).Value = '  +1.58%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -0.02%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4841'
$ws.Cells.Item(7, 5).Value = '  +3.46%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +1.68%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06730'
$ws.Cells.Item(9, 5).Value = '  -1.03%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '110.12'
$ws.Cells.Item(10, 5).Value = '  +2.56%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '19.04'
$ws.Cells.Item(11, 5).Value = '  +4.87%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.918.62'
$ws.Cells.Item(12, 5).Value = '  +0.09%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.07551'
$ws.Cells.Item(13, 5).Value = '  -1.10%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.277'
$ws.Cells.Item(14, 5).Value = '  +1.95%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.6726'
$ws.Cells.Item(15, 5).Value = '  +2.70%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '280.25'
$ws.Cells.Item(16, 5).Value = '  -2.78%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '30.485.12'
$ws.Cells.Item(17, 5).Value = '  +0.22%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '1.000'
$ws.Cells.Item(18, 5).Value = '  +0.08%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.000007570'
$ws.Cells.Item(19, 5).Value = '  -0.45%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.87'
$ws.Cells.Item(20, 5).Value = '  -0.93%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +5.66%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '2.165.19'
$ws.Cells.Item(22, 5).Value = '  +0.30%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.0000'
$ws.Cells.Item(23, 5).Value = '  +0.05%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.466'
$ws.Cells.Item(24, 5).Value = '  +4.39%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.453'
$ws.Cells.Item(25, 5).Value = '  +2.04%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '163.79'
$ws.Cells.Item(26, 5).Value = '  -2.57%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '20.21'
$ws.Cells.Item(27, 5).Value = '  -6.56%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.128'
$ws.Cells.Item(28, 5).Value = '  +3.76%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.1056'
$ws.Cells.Item(29, 5).Value = '  -1.31%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.401'
$ws.Cells.Item(30, 5).Value = '  +2.03%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.152'
$ws.Cells.Item(31, 5).Value = '  +0.45%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.047'
$ws.Cells.Item(32, 5).Value = '  +2.68%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.04992'
$ws.Cells.Item(33, 5).Value = '  -0.83%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.7316'
$ws.Cells.Item(34, 5).Value = '  -0.82%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.135'
$ws.Cells.Item(35, 5).Value = '  -0.96%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.9995'
$ws.Cells.Item(36, 5).Value = '  +0.10%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.729'
$ws.Cells.Item(37, 5).Value = '  -0.16%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.21%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.666'
$ws.Cells.Item(39, 5).Value = '  -0.61%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '111.11'
$ws.Cells.Item(40, 5).Value = '  +2.10%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.015'
$ws.Cells.Item(41, 5).Value = '  -1.80%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.4442'
$ws.Cells.Item(42, 5).Value = '  +5.70%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.8657'
$ws.Cells.Item(43, 5).Value = '  -0.84%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '5.799'
$ws.Cells.Item(44, 5).Value = '  -0.38%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.9997'
$ws.Cells.Item(45, 5).Value = '  +0.01%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '68.09'
$ws.Cells.Item(46, 5).Value = '  +0.86%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '7.349'
$ws.Cells.Item(47, 5).Value = '  +2.63%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '9.247'
$ws.Cells.Item(48, 5).Value = '  +0.19%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '48.06'
$ws.Cells.Item(49, 5).Value = '  -9.17%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.1242'
$ws.Cells.Item(50, 5).Value = '  +2.91%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.469'
$ws.Cells.Item(51, 5).Value = '  +5.90%  '
